$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'318.16"
$ws.Range("E2").Value = "'3.08%"
$ws.Range("D3").Value = "'41.46"
$ws.Range("E3").Value = "'1.36%"
$ws.Range("E4").Value = "'2.41%"
$ws.Range("D5").Value = "'0.07734"
$ws.Range("D6").Value = "'1.693"
$ws.Range("E6").Value = "'4.31%"
$ws.Range("D7").Value = "'0.9519"
$ws.Range("E7").Value = "'4.66%"
$ws.Range("E8").Value = "'-1.30%"
$ws.Range("D9").Value = "'0.1261"
$ws.Range("E9").Value = "'6.66%"
$ws.Range("D10").Value = "'0.1836"
$ws.Range("E10").Value = "'0.87%"
$ws.Range("D11").Value = "'0.09179"
$ws.Range("E11").Value = "'-0.25%"
$ws.Range("D12").Value = "'0.04411"
$ws.Range("E12").Value = "'3.54%"
$ws.Range("E13").Value = "'0.54%"
$ws.Range("D14").Value = "'0.001285"
$ws.Range("E14").Value = "'2.85%"
$ws.Range("D15").Value = "'0.006000"
$ws.Range("E15").Value = "'3.50%"
$ws.Range("D16").Value = "'3.337"
$ws.Range("E16").Value = "'-0.42%"
$ws.Range("D17").Value = "'4.325"
$ws.Range("E17").Value = "'1.03%"
$ws.Range("D19").Value = "'7.664"
$ws.Range("E19").Value = "'10.93%"
$ws.Range("D20").Value = "'0.1352"
$ws.Range("E20").Value = "'-4.17%"
$ws.Range("D21").Value = "'0.2822"
$ws.Range("E21").Value = "'4.30%"
$ws.Range("D22").Value = "'0.04023"
$ws.Range("E22").Value = "'-0.35%"
$ws.Range("D23").Value = "'0.001264"
$ws.Range("E23").Value = "'-0.58%"
$ws.Range("D24").Value = "'0.004125"
$ws.Range("E24").Value = "'0.41%"
$ws.Range("E25").Value = "'-0.24%"
$ws.Range("D38").Value = "'0.02547"
$ws.Range("E38").Value = "'5.07%"
$ws.Range("D39").Value = "'0.05352"
$ws.Range("E39").Value = "'2.10%"
$ws.Range("D40").Value = "'0.007778"
$ws.Range("E40").Value = "'-0.27%"
$ws.Range("D41").Value = "'0.1319"
$ws.Range("E41").Value = "'1.43%"
$ws.Range("D42").Value = "'0.007329"
$ws.Range("E42").Value = "'7.69%"
$ws.Range("D43").Value = "'0.001990"
$ws.Range("E43").Value = "'2.90%"
$ws.Range("D44").Value = "'0.007566"
$ws.Range("E44").Value = "'-6.36%"
$ws.Range("D45").Value = "'0.3433"
$ws.Range("E45").Value = "'12.01%"
$ws.Range("D46").Value = "'0.00006686"
$ws.Range("E46").Value = "'-3.08%"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("D48").Value = "'0.1933"
$ws.Range("E48").Value = "'86.07%"
$ws.Range("D49").Value = "'0.004200"
$ws.Range("E49").Value = "'39.88%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.20%"
